$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure ID column (D) retains leading zeros as text, matching the
# original "t=s" (string) cell type used throughout the template.
$idRange = $ws.Range("D2:D31")
$idRange.NumberFormat = "@"

$ws.Range("A2").Value = "Proctor44143"
$ws.Range("B2").Value = "Automation44143"
$ws.Range("C2").Value = "proctorautomation44143@gmail.com"
$ws.Range("D2").Value = "44143"
$ws.Range("A3").Value = "Proctor09201"
$ws.Range("B3").Value = "Automation09201"
$ws.Range("C3").Value = "proctorautomation09201@gmail.com"
$ws.Range("D3").Value = "09201"
$ws.Range("A4").Value = "Proctor29544"
$ws.Range("B4").Value = "Automation29544"
$ws.Range("C4").Value = "proctorautomation29544@gmail.com"
$ws.Range("D4").Value = "29544"
$ws.Range("A5").Value = "Proctor78719"
$ws.Range("B5").Value = "Automation78719"
$ws.Range("C5").Value = "proctorautomation78719@gmail.com"
$ws.Range("D5").Value = "78719"
$ws.Range("A6").Value = "Proctor85406"
$ws.Range("B6").Value = "Automation85406"
$ws.Range("C6").Value = "proctorautomation85406@gmail.com"
$ws.Range("D6").Value = "85406"
$ws.Range("A7").Value = "Proctor43828"
$ws.Range("B7").Value = "Automation43828"
$ws.Range("C7").Value = "proctorautomation43828@gmail.com"
$ws.Range("D7").Value = "43828"
$ws.Range("A8").Value = "Proctor76382"
$ws.Range("B8").Value = "Automation76382"
$ws.Range("C8").Value = "proctorautomation76382@gmail.com"
$ws.Range("D8").Value = "76382"
$ws.Range("A9").Value = "Proctor43694"
$ws.Range("B9").Value = "Automation43694"
$ws.Range("C9").Value = "proctorautomation43694@gmail.com"
$ws.Range("D9").Value = "43694"
$ws.Range("A10").Value = "Proctor12498"
$ws.Range("B10").Value = "Automation12498"
$ws.Range("C10").Value = "proctorautomation12498@gmail.com"
$ws.Range("D10").Value = "12498"
$ws.Range("A11").Value = "Proctor93448"
$ws.Range("B11").Value = "Automation93448"
$ws.Range("C11").Value = "proctorautomation93448@gmail.com"
$ws.Range("D11").Value = "93448"
$ws.Range("A12").Value = "Proctor21145"
$ws.Range("B12").Value = "Automation21145"
$ws.Range("C12").Value = "proctorautomation21145@gmail.com"
$ws.Range("D12").Value = "21145"
$ws.Range("A13").Value = "Proctor22176"
$ws.Range("B13").Value = "Automation22176"
$ws.Range("C13").Value = "proctorautomation22176@gmail.com"
$ws.Range("D13").Value = "22176"
$ws.Range("A14").Value = "Proctor13242"
$ws.Range("B14").Value = "Automation13242"
$ws.Range("C14").Value = "proctorautomation13242@gmail.com"
$ws.Range("D14").Value = "13242"
$ws.Range("A15").Value = "Proctor18017"
$ws.Range("B15").Value = "Automation18017"
$ws.Range("C15").Value = "proctorautomation18017@gmail.com"
$ws.Range("D15").Value = "18017"
$ws.Range("A16").Value = "Proctor62907"
$ws.Range("B16").Value = "Automation62907"
$ws.Range("C16").Value = "proctorautomation62907@gmail.com"
$ws.Range("D16").Value = "62907"
$ws.Range("A17").Value = "Proctor83774"
$ws.Range("B17").Value = "Automation83774"
$ws.Range("C17").Value = "proctorautomation83774@gmail.com"
$ws.Range("D17").Value = "83774"
$ws.Range("A18").Value = "Proctor70984"
$ws.Range("B18").Value = "Automation70984"
$ws.Range("C18").Value = "proctorautomation70984@gmail.com"
$ws.Range("D18").Value = "70984"
$ws.Range("A19").Value = "Proctor52186"
$ws.Range("B19").Value = "Automation52186"
$ws.Range("C19").Value = "proctorautomation52186@gmail.com"
$ws.Range("D19").Value = "52186"
$ws.Range("A20").Value = "Proctor07467"
$ws.Range("B20").Value = "Automation07467"
$ws.Range("C20").Value = "proctorautomation07467@gmail.com"
$ws.Range("D20").Value = "07467"
$ws.Range("A21").Value = "Proctor49207"
$ws.Range("B21").Value = "Automation49207"
$ws.Range("C21").Value = "proctorautomation49207@gmail.com"
$ws.Range("D21").Value = "49207"
$ws.Range("A22").Value = "Proctor74337"
$ws.Range("B22").Value = "Automation74337"
$ws.Range("C22").Value = "proctorautomation74337@gmail.com"
$ws.Range("D22").Value = "74337"
$ws.Range("A23").Value = "Proctor11708"
$ws.Range("B23").Value = "Automation11708"
$ws.Range("C23").Value = "proctorautomation11708@gmail.com"
$ws.Range("D23").Value = "11708"
$ws.Range("A24").Value = "Proctor04159"
$ws.Range("B24").Value = "Automation04159"
$ws.Range("C24").Value = "proctorautomation04159@gmail.com"
$ws.Range("D24").Value = "04159"
$ws.Range("A25").Value = "Proctor19308"
$ws.Range("B25").Value = "Automation19308"
$ws.Range("C25").Value = "proctorautomation19308@gmail.com"
$ws.Range("D25").Value = "19308"
$ws.Range("A26").Value = "Proctor73696"
$ws.Range("B26").Value = "Automation73696"
$ws.Range("C26").Value = "proctorautomation73696@gmail.com"
$ws.Range("D26").Value = "73696"
$ws.Range("A27").Value = "Proctor07711"
$ws.Range("B27").Value = "Automation07711"
$ws.Range("C27").Value = "proctorautomation07711@gmail.com"
$ws.Range("D27").Value = "07711"
$ws.Range("A28").Value = "Proctor32355"
$ws.Range("B28").Value = "Automation32355"
$ws.Range("C28").Value = "proctorautomation32355@gmail.com"
$ws.Range("D28").Value = "32355"
$ws.Range("A29").Value = "Proctor63231"
$ws.Range("B29").Value = "Automation63231"
$ws.Range("C29").Value = "proctorautomation63231@gmail.com"
$ws.Range("D29").Value = "63231"
$ws.Range("A30").Value = "Proctor81924"
$ws.Range("B30").Value = "Automation81924"
$ws.Range("C30").Value = "proctorautomation81924@gmail.com"
$ws.Range("D30").Value = "81924"
$ws.Range("A31").Value = "Proctor74585"
$ws.Range("B31").Value = "Automation74585"
$ws.Range("C31").Value = "proctorautomation74585@gmail.com"
$ws.Range("D31").Value = "74585"

# Drop the explicit text-format style so the cells match the
# workbook default (unstyled) formatting used by the other data rows.
$idRange.ClearFormats()
